$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 8, pushing nothing (rows 6/7 currently end the data at row 7).
# This creates blank rows 8 and 9 below current row 7, preserving rows 1-7 as-is.
$ws.Rows("8:9").Insert()

# Populate new rows 8 & 9 with the data that used to live in rows 6 & 7
# (the weekly snapshot being superseded), each cell set explicitly for
# robustness instead of relying on clipboard copy/paste semantics.

# Row 8 = old Row 6 data (Especial, volumen 45, price 14000, 7 kilos)
$ws.Cells.Item(8, 1).Value = 9
$ws.Cells.Item(8, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44187
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100101
$ws.Cells.Item(8, 8).Value = "Berries"
$ws.Cells.Item(8, 9).Value = 100101006
$ws.Cells.Item(8, 10).Value = "Breva"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Especial"
$ws.Cells.Item(8, 13).Value = 45
$ws.Cells.Item(8, 14).Value = 14000
$ws.Cells.Item(8, 15).Value = 14000
$ws.Cells.Item(8, 16).Value = 14000
$ws.Cells.Item(8, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(8, 19).Value = 2000
$ws.Cells.Item(8, 20).Value = 7

# Row 9 = old Row 7 data (Primera, volumen 50, price 12000, 7 kilos)
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44187
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100101
$ws.Cells.Item(9, 8).Value = "Berries"
$ws.Cells.Item(9, 9).Value = 100101006
$ws.Cells.Item(9, 10).Value = "Breva"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 50
$ws.Cells.Item(9, 14).Value = 12000
$ws.Cells.Item(9, 15).Value = 12000
$ws.Cells.Item(9, 16).Value = 12000
$ws.Cells.Item(9, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(9, 19).Value = 1714
$ws.Cells.Item(9, 20).Value = 7

# Now update row 6 with the new values (newer date, new volume/prices, 6 kilos)
$ws.Cells.Item(6, 4).Value = 44553
$ws.Cells.Item(6, 13).Value = 200
$ws.Cells.Item(6, 14).Value = 22000
$ws.Cells.Item(6, 15).Value = 22000
$ws.Cells.Item(6, 16).Value = 22000
$ws.Cells.Item(6, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(6, 19).Value = 3667
$ws.Cells.Item(6, 20).Value = 6

# Update row 7 with the new values
$ws.Cells.Item(7, 4).Value = 44553
$ws.Cells.Item(7, 13).Value = 150
$ws.Cells.Item(7, 14).Value = 18000
$ws.Cells.Item(7, 15).Value = 18000
$ws.Cells.Item(7, 16).Value = 18000
$ws.Cells.Item(7, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(7, 19).Value = 3000
$ws.Cells.Item(7, 20).Value = 6
